# Update "想去人数" (number of people interested) figures for a handful of
# events, on both the "展览" sheet and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 271
$ws1.Range("F4").Value = 1932
$ws1.Range("F5").Value = 1585
$ws1.Range("F6").Value = 286
$ws1.Range("F8").Value = 587
$ws1.Range("F9").Value = 135

# --- Sheet "全部类型" (all types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 271
$ws4.Range("F4").Value = 1932
$ws4.Range("F5").Value = 1585
$ws4.Range("F6").Value = 286
$ws4.Range("F9").Value = 587
$ws4.Range("F10").Value = 135
